$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Main"
$ws.Range("B1").Value = "Xpath"
$ws.Range("C1").Value = "Value"

# Data rows (Xpath / Value pairs)
$ws.Range("B2").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[1]/CITY"
$ws.Range("C2").Value = "[A-Z a-z].*"

$ws.Range("B3").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[1]/STATE"
$ws.Range("C3").Value = "[A-Z]{2}"

$ws.Range("B4").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[1]/ZIP"
$ws.Range("C4").Value = "[0-9]{5}"

$ws.Range("B5").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[1]/AREA_CODE"
$ws.Range("C5").Value = "[0-9]{3}"

$ws.Range("B6").Value = "/Envelope/Body/GetInfoByAreaCodeResponse/GetInfoByAreaCodeResult/NewDataSet/Table[1]/TIME_ZONE"
$ws.Range("C6").Value = "[A-Z]{1}"

# Column widths (Excel's stored XML width = ColumnWidth + 0.8333, so back it out)
$ws.Columns.Item(2).ColumnWidth = 100.16666666666667
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666

# Selection matching the final view state
$ws.Range("C2:C6").Select()
